# Case_4_29 (380 kV) results refresh: res_bus/vm_pu.xlsx
# Slack-bus voltage setpoint was changed from 1.05 to 1.02 p.u., so every
# per-unit voltage magnitude result in rows 2-25 (columns B-F, I-N) is updated
# to the corresponding recomputed power-flow value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.080545938900384
$ws.Range("D2").Value = 1.080358253600873
$ws.Range("E2").Value = 1.084281696624755
$ws.Range("F2").Value = 1.090388847334099
$ws.Range("I2").Value = 1.066517180749263
$ws.Range("J2").Value = 1.085424939894193
$ws.Range("K2").Value = 1.083031896845428
$ws.Range("L2").Value = 1.086945122710575
$ws.Range("M2").Value = 1.093036527425253
$ws.Range("N2").Value = 1.086966367808823

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.082076833778985
$ws.Range("D3").Value = 1.081607143058637
$ws.Range("E3").Value = 1.085765685653733
$ws.Range("F3").Value = 1.091779194970309
$ws.Range("I3").Value = 1.067138822253866
$ws.Range("J3").Value = 1.086614527559414
$ws.Range("K3").Value = 1.084097984020297
$ws.Range("L3").Value = 1.088246486689164
$ws.Range("M3").Value = 1.094245622910291
$ws.Range("N3").Value = 1.088157644824977

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.083064632568904
$ws.Range("D4").Value = 1.082412283707002
$ws.Range("E4").Value = 1.086723453255276
$ws.Range("F4").Value = 1.092676262388592
$ws.Range("I4").Value = 1.067538137897128
$ws.Range("J4").Value = 1.087381031350491
$ws.Range("K4").Value = 1.084784277759929
$ws.Range("L4").Value = 1.089085566952537
$ws.Range("M4").Value = 1.095024870617743
$ws.Range("N4").Value = 1.088925237139357

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.08347924587696
$ws.Range("D5").Value = 1.082750061916173
$ws.Range("E5").Value = 1.087125517166859
$ws.Range("F5").Value = 1.093052781405889
$ws.Range("I5").Value = 1.067705314691507
$ws.Range("J5").Value = 1.087702503318073
$ws.Range("K5").Value = 1.0850719578071
$ws.Range("L5").Value = 1.089437610318678
$ws.Range("M5").Value = 1.095351728993833
$ws.Range("N5").Value = 1.089247165634006

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.083548823056105
$ws.Range("D6").Value = 1.082806735408443
$ws.Range("E6").Value = 1.087192991661007
$ws.Range("F6").Value = 1.093115965223578
$ws.Range("I6").Value = 1.067733343784087
$ws.Range("J6").Value = 1.087756435209645
$ws.Range("K6").Value = 1.085120211672037
$ws.Range("L6").Value = 1.089496678829118
$ws.Range("M6").Value = 1.095406566976512
$ws.Range("N6").Value = 1.08930117411505

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.083070175221079
$ws.Range("D7").Value = 1.082416799866656
$ws.Range("E7").Value = 1.086728827925402
$ws.Range("F7").Value = 1.092681295831722
$ws.Range("I7").Value = 1.067540374447753
$ws.Range("J7").Value = 1.087385329874055
$ws.Range("K7").Value = 1.084788125036231
$ws.Range("L7").Value = 1.089090273735305
$ws.Range("M7").Value = 1.095029241001471
$ws.Range("N7").Value = 1.088929541767316

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.08106389700602
$ws.Range("D8").Value = 1.080780943025588
$ws.Range("E8").Value = 1.084783735869806
$ws.Range("F8").Value = 1.090859261961709
$ws.Range("I8").Value = 1.066727876965933
$ws.Range("J8").Value = 1.085827642822208
$ws.Range("K8").Value = 1.083392923382314
$ws.Range("L8").Value = 1.08738554954124
$ws.Range("M8").Value = 1.093445797978768
$ws.Range("N8").Value = 1.087369642621187

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.077506684000913
$ws.Range("D9").Value = 1.077875161882819
$ws.Range("E9").Value = 1.08133683225285
$ws.Range("F9").Value = 1.087628407148524
$ws.Range("I9").Value = 1.065273481304203
$ws.Range("J9").Value = 1.08305757803457
$ws.Range("K9").Value = 1.080906927472517
$ws.Range("L9").Value = 1.084358278162774
$ws.Range("M9").Value = 1.090631280101038
$ws.Range("N9").Value = 1.084595644023819

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.075119775606769
$ws.Range("D10").Value = 1.075921800623693
$ws.Range("E10").Value = 1.079025195201105
$ws.Range("F10").Value = 1.085460311506787
$ws.Range("I10").Value = 1.064288294207324
$ws.Range("J10").Value = 1.08119332761704
$ws.Range("K10").Value = 1.079230575522828
$ws.Range("L10").Value = 1.082323799109974
$ws.Range("M10").Value = 1.088738028121953
$ws.Range("N10").Value = 1.082728746156864

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.074082397761838
$ws.Range("D11").Value = 1.075072007123075
$ws.Range("E11").Value = 1.078020833995269
$ws.Range("F11").Value = 1.08451799853815
$ws.Range("I11").Value = 1.063857921625797
$ws.Range("J11").Value = 1.080381795378279
$ws.Range("K11").Value = 1.078500061501948
$ws.Range("L11").Value = 1.081438844103074
$ws.Range("M11").Value = 1.087914089770457
$ws.Range("N11").Value = 1.08191606144924

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.073696480049059
$ws.Range("D12").Value = 1.074755747413307
$ws.Range("E12").Value = 1.077647244863686
$ws.Range("F12").Value = 1.084167441771554
$ws.Range("I12").Value = 1.063697487371939
$ws.Range("J12").Value = 1.080079698834077
$ws.Range("K12").Value = 1.078228007821489
$ws.Range("L12").Value = 1.081109517678963
$ws.Range("M12").Value = 1.087607407975051
$ws.Range("N12").Value = 1.081613535893295

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.073779287603901
$ws.Range("D13").Value = 1.07482361391697
$ws.Range("E13").Value = 1.077727404949298
$ws.Range("F13").Value = 1.084242662007126
$ws.Range("I13").Value = 1.063731927181081
$ws.Range("J13").Value = 1.080144529498337
$ws.Range("K13").Value = 1.078286396483061
$ws.Range("L13").Value = 1.081180187326075
$ws.Range("M13").Value = 1.087673221151436
$ws.Range("N13").Value = 1.081678458624534

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.074050509787462
$ws.Range("D14").Value = 1.075045877477432
$ws.Range("E14").Value = 1.077989963778463
$ws.Range("F14").Value = 1.084489032487271
$ws.Range("I14").Value = 1.063844671843954
$ws.Range("J14").Value = 1.080356837459486
$ws.Range("K14").Value = 1.078477587970511
$ws.Range("L14").Value = 1.081411634488075
$ws.Range("M14").Value = 1.087888752364808
$ws.Range("N14").Value = 1.081891068087339

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.074217540208759
$ws.Range("D15").Value = 1.075182740410272
$ws.Range("E15").Value = 1.078151665031278
$ws.Range("F15").Value = 1.084640757586414
$ws.Range("I15").Value = 1.063914061191712
$ws.Range("J15").Value = 1.080487559955085
$ws.Range("K15").Value = 1.078595293182697
$ws.Range("L15").Value = 1.081554154905566
$ws.Range("M15").Value = 1.088021463855086
$ws.Range("N15").Value = 1.082021976223877

# Row 16
$ws.Range("B16").Value = 1.019999999999999
$ws.Range("C16").Value = 1.075188541068201
$ws.Range("D16").Value = 1.075978113901977
$ws.Range("E16").Value = 1.079091778474964
$ws.Range("F16").Value = 1.085522774678199
$ws.Range("I16").Value = 1.064316776397606
$ws.Range("J16").Value = 1.081247094784449
$ws.Range("K16").Value = 1.079278958603792
$ws.Range("L16").Value = 1.082382445110706
$ws.Range("M16").Value = 1.088792621850755
$ws.Range("N16").Value = 1.082782589679819

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.075796588910236
$ws.Range("D17").Value = 1.076475958443659
$ws.Range("E17").Value = 1.079680565641801
$ws.Range("F17").Value = 1.086075091707543
$ws.Range("I17").Value = 1.064568372196805
$ws.Range("J17").Value = 1.081722371818151
$ws.Range("K17").Value = 1.079706553358893
$ws.Range("L17").Value = 1.082900927110012
$ws.Range("M17").Value = 1.089275230511596
$ws.Range("N17").Value = 1.083258541661423

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.076150884260745
$ws.Range("D18").Value = 1.076765959966979
$ws.Range("E18").Value = 1.080023667368995
$ws.Range("F18").Value = 1.086396910722902
$ws.Range("I18").Value = 1.064714759511688
$ws.Range("J18").Value = 1.081999178806597
$ws.Range("K18").Value = 1.079955514937777
$ws.Range("L18").Value = 1.083202962588025
$ws.Range("M18").Value = 1.089556328547328
$ws.Range("N18").Value = 1.083535741747543

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.076271627703991
$ws.Range("D19").Value = 1.076864778407905
$ws.Range("E19").Value = 1.080140600964164
$ws.Range("F19").Value = 1.086506585691556
$ws.Range("I19").Value = 1.06476461222316
$ws.Range("J19").Value = 1.082093492946299
$ws.Range("K19").Value = 1.080040328871774
$ws.Range("L19").Value = 1.083305883717018
$ws.Range("M19").Value = 1.089652108301319
$ws.Range("N19").Value = 1.083630189824141

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.075731389364543
$ws.Range("D20").Value = 1.076422584110023
$ws.Range("E20").Value = 1.079617428356054
$ws.Range("F20").Value = 1.086015868381516
$ws.Range("I20").Value = 1.064541416074692
$ws.Range("J20").Value = 1.081671422001261
$ws.Range("K20").Value = 1.079660722822831
$ws.Range("L20").Value = 1.082845338941526
$ws.Range("M20").Value = 1.089223492565969
$ws.Range("N20").Value = 1.08320751948995

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.073970657996625
$ws.Range("D21").Value = 1.074980443289766
$ws.Range("E21").Value = 1.077912661333923
$ws.Range("F21").Value = 1.084416497518496
$ws.Range("I21").Value = 1.063811487270773
$ws.Range("J21").Value = 1.080294336294348
$ws.Range("K21").Value = 1.07842130648082
$ws.Range("L21").Value = 1.081343496133782
$ws.Range("M21").Value = 1.087825301398144
$ws.Range("N21").Value = 1.081828478163379

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.072860199286828
$ws.Range("D22").Value = 1.074070185430475
$ws.Range("E22").Value = 1.076837764688488
$ws.Range("F22").Value = 1.083407780600545
$ws.Range("I22").Value = 1.063349223514169
$ws.Range("J22").Value = 1.079424698972959
$ws.Range("K22").Value = 1.077637933122898
$ws.Range("L22").Value = 1.080395665287909
$ws.Range("M22").Value = 1.086942526976053
$ws.Range("N22").Value = 1.080957605857244

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.073449203020026
$ws.Range("D23").Value = 1.074553068451896
$ws.Range("E23").Value = 1.077407880383323
$ws.Range("F23").Value = 1.083942821138065
$ws.Range("I23").Value = 1.063594596076184
$ws.Range("J23").Value = 1.079886075042297
$ws.Range("K23").Value = 1.078053606710242
$ws.Range("L23").Value = 1.080898470076193
$ws.Range("M23").Value = 1.087410854644503
$ws.Range("N23").Value = 1.081419637133521

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.075760851371869
$ws.Range("D24").Value = 1.076446702856657
$ws.Range("E24").Value = 1.079645958390843
$ws.Range("F24").Value = 1.086042629898215
$ws.Range("I24").Value = 1.064553597511464
$ws.Range("J24").Value = 1.081694445309796
$ws.Range("K24").Value = 1.079681433051117
$ws.Range("L24").Value = 1.082870458035712
$ws.Range("M24").Value = 1.089246871950756
$ws.Range("N24").Value = 1.083230575494223

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.078428976330021
$ws.Range("D25").Value = 1.078629183760722
$ws.Range("E25").Value = 1.082230305902759
$ws.Range("F25").Value = 1.088466116791595
$ws.Range("I25").Value = 1.065652199607555
$ws.Range("J25").Value = 1.083776754683178
$ws.Range("K25").Value = 1.081552928927151
$ws.Range("L25").Value = 1.085143727721681
$ws.Range("M25").Value = 1.091361837741381
$ws.Range("N25").Value = 1.08531584198576
